# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.273.35"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.594.29"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "1.819.01"
$ws.Range("D13").Value = "1.612.16"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.267.85"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.33%  "
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("B32").Value = "Maker"
$ws.Range("C32").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D32").Value = "1.471.94"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.815"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.926"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "1.731.61"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  +0.00%  "
